$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The old row 6 (S022856 / party 542161) reported a single "water year"
# 2022 entry from one measurement spreadsheet. The QAQC script now
# corrects calendar-year and water-year rows differently, so that one
# row is split in two: a 2021 calendar-year row (new row 6, only the
# OCT/NOV/DEC columns O/P/Q apply) and the original 2022 row (pushed
# down to row 7, unchanged).
$ws.Rows(6).Insert()

# New row 6: 2021 portion of the same right/spreadsheet.
$ws.Range("A6").Value = "S022856"
$ws.Range("B6").Value = 542161
$ws.Range("C6").Value = 2021
$ws.Range("D6").Value = "https://rms.waterboards.ca.gov/DownloadDocument.ashx?type=Attachment&download=true&id=85316"
$ws.Range("E6").Value = "DIRECT"
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0

# Same source measurement spreadsheet as row 7 - link it with its own
# hyperlink relationship (Excel gives every row its own rId even when
# two rows point at the same URL, as already done for rows 4 & 5).
$ws.Hyperlinks.Add($ws.Range("D6"), "https://rms.waterboards.ca.gov/DownloadDocument.ashx?type=Attachment&download=true&id=85316") | Out-Null
$ws.Range("D6").Style = "Hyperlink"
